# Auto-generated Excel COM-interop script to apply the Gungnir_Profits diff
$wb = $excel.ActiveWorkbook

# hunk 0  sheet=ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 325
$ws.Range("I61").Value = 266.66666
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 799.9999799999999
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -627.9999799999999
$ws.Range("N61").Value = -1844

# hunk 1  sheet=ALC
$ws.Range("H64").Value = 3541.6667
$ws.Range("I64").Value = 3214.2856
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3214.2856
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -2966.2856
$ws.Range("N64").Value = -4496

# hunk 2  sheet=ALC
$ws.Range("H67").Value = 3541.6667
$ws.Range("I67").Value = 3214.2856
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3214.2856
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2356.2856
$ws.Range("N67").Value = -5716

# hunk 3  sheet=ALC
$ws.Range("H76").Value = 10107044
$ws.Range("I76").Value = 10091.429
$ws.Range("K76").Value = 10091.429
$ws.Range("M76").Value = -9776.429

# hunk 4  sheet=ALC
$ws.Range("H79").Value = 10107044
$ws.Range("I79").Value = 10091.429
$ws.Range("K79").Value = 10091.429
$ws.Range("M79").Value = -8999.429

# hunk 5  sheet=ALC
$ws.Range("H86").Value = 15127.219
$ws.Range("I86").Value = 12611.4
$ws.Range("J86").Value = 17347.059
$ws.Range("K86").Value = 12611.4
$ws.Range("L86").Value = 17347.059
$ws.Range("M86").Value = -11488.4
$ws.Range("N86").Value = -19593.059

# hunk 6  sheet=ALC
$ws.Range("H89").Value = 15127.219
$ws.Range("I89").Value = 12611.4
$ws.Range("J89").Value = 17347.059
$ws.Range("K89").Value = 63057
$ws.Range("L89").Value = 86735.29500000001
$ws.Range("M89").Value = -57441
$ws.Range("N89").Value = -97967.29500000001

# hunk 7  sheet=ALC
$ws.Range("H92").Value = 2525655.8
$ws.Range("I92").Value = 2778171.5
$ws.Range("K92").Value = 2778171.5
$ws.Range("M92").Value = -2776923.5

# hunk 8  sheet=ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13515702
$ws.Range("I32").Value = 1935.2206
$ws.Range("J32").Value = 166671730
$ws.Range("K32").Value = 1935.2206
$ws.Range("L32").Value = 166671730
$ws.Range("M32").Value = -1648.2206
$ws.Range("N32").Value = -166672304

# hunk 9  sheet=ARM
$ws.Range("H44").Value = 47780
$ws.Range("J44").Value = 47780
$ws.Range("L44").Value = 47780
$ws.Range("N44").Value = -48756

# hunk 10  sheet=ARM
$ws.Range("H63").Value = 2851.121
$ws.Range("I63").Value = 2925.423
$ws.Range("K63").Value = 2925.423
$ws.Range("M63").Value = -2239.423

# hunk 11  sheet=ARM
$ws.Range("H66").Value = 2851.121
$ws.Range("I66").Value = 2925.423
$ws.Range("K66").Value = 14627.115
$ws.Range("M66").Value = -11195.115

# hunk 12  sheet=ARM
$ws.Range("H132").Value = 1839645.8
$ws.Range("I132").Value = 1064.8846
$ws.Range("J132").Value = 9806829
$ws.Range("K132").Value = 3194.6538
$ws.Range("L132").Value = 29420487
$ws.Range("M132").Value = -664.6538
$ws.Range("N132").Value = -29425547

# hunk 13  sheet=BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1702.9482
$ws.Range("I20").Value = 1649.0264
$ws.Range("J20").Value = 1805.4
$ws.Range("K20").Value = 1649.0264
$ws.Range("L20").Value = 1805.4
$ws.Range("M20").Value = -1402.0264
$ws.Range("N20").Value = -2299.4

# hunk 14  sheet=BSM
$ws.Range("H105").Value = 62501936
$ws.Range("I105").Value = 1907.6666
$ws.Range("J105").Value = 142859120
$ws.Range("K105").Value = 1907.6666
$ws.Range("L105").Value = 142859120
$ws.Range("M105").Value = -160.6666
$ws.Range("N105").Value = -142862614

# hunk 15  sheet=BSM
$ws.Range("H134").Value = 1483782.6
$ws.Range("I134").Value = 909.50793
$ws.Range("J134").Value = 9268866
$ws.Range("K134").Value = 2728.52379
$ws.Range("L134").Value = 27806598
$ws.Range("M134").Value = -193.5237900000002
$ws.Range("N134").Value = -27811668

# hunk 16  sheet=CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9574.25
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 11099
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 11099
$ws.Range("M51").Value = -4264
$ws.Range("N51").Value = -12571

# hunk 17  sheet=CRP
$ws.Range("H61").Value = 9574.25
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 11099
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 11099
$ws.Range("M61").Value = -4652
$ws.Range("N61").Value = -11795

# hunk 18  sheet=CRP
$ws.Range("H134").Value = 936.1613
$ws.Range("I134").Value = 1060.9048
$ws.Range("J134").Value = 674.2
$ws.Range("K134").Value = 3182.7144
$ws.Range("L134").Value = 2022.6
$ws.Range("M134").Value = -647.7143999999998
$ws.Range("N134").Value = -7092.6

# hunk 19  sheet=CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2665
$ws.Range("I55").Value = 800
$ws.Range("J55").Value = 3038
$ws.Range("K55").Value = 2400
$ws.Range("L55").Value = 9114
$ws.Range("M55").Value = -2223
$ws.Range("N55").Value = -9468

# hunk 20  sheet=CUL
$ws.Range("H113").Value = 2326094.2
$ws.Range("I113").Value = 535.44446
$ws.Range("J113").Value = 4000496.8
$ws.Range("K113").Value = 1606.33338
$ws.Range("L113").Value = 12001490.4
$ws.Range("M113").Value = 563.66662
$ws.Range("N113").Value = -12005830.4

# hunk 21  sheet=CUL
$ws.Range("H115").Value = 3898.4167
$ws.Range("I115").Value = 447
$ws.Range("J115").Value = 5624.125
$ws.Range("K115").Value = 1341
$ws.Range("L115").Value = 16872.375
$ws.Range("M115").Value = -166
$ws.Range("N115").Value = -19222.375

# hunk 22  sheet=CUL
$ws.Range("H131").Value = 879.9
$ws.Range("I131").Value = 665
$ws.Range("J131").Value = 884.2857
$ws.Range("K131").Value = 1995
$ws.Range("L131").Value = 2652.8571
$ws.Range("M131").Value = 3045
$ws.Range("N131").Value = -12732.8571

# hunk 23  sheet=CUL
$ws.Range("H137").Value = 29413296
$ws.Range("I137").Value = 45455316
$ws.Range("J137").Value = 2927.5
$ws.Range("K137").Value = 136365948
$ws.Range("L137").Value = 8782.5
$ws.Range("M137").Value = -136360848
$ws.Range("N137").Value = -18982.5

# hunk 24  sheet=GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# hunk 25  sheet=GSM
$ws.Range("H70").Value = 11211.5
$ws.Range("I70").Value = 15106.777
$ws.Range("J70").Value = 4200
$ws.Range("K70").Value = 15106.777
$ws.Range("L70").Value = 4200
$ws.Range("M70").Value = -14836.777
$ws.Range("N70").Value = -4740

# hunk 26  sheet=GSM
$ws.Range("H73").Value = 11211.5
$ws.Range("I73").Value = 15106.777
$ws.Range("J73").Value = 4200
$ws.Range("K73").Value = 15106.777
$ws.Range("L73").Value = 4200
$ws.Range("M73").Value = -14170.777
$ws.Range("N73").Value = -6072

# hunk 27  sheet=LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 178572670
$ws.Range("I16").Value = 17857754
$ws.Range("K16").Value = 17857754
$ws.Range("M16").Value = -17857584

# hunk 28  sheet=LTW
$ws.Range("H17").Value = 6335.2
$ws.Range("I17").Value = 2919
$ws.Range("K17").Value = 2919
$ws.Range("M17").Value = -2749

# hunk 29  sheet=LTW
$ws.Range("H22").Value = 975996.1
$ws.Range("I22").Value = 2533901.5
$ws.Range("K22").Value = 2533901.5
$ws.Range("M22").Value = -2533606.5

# hunk 30  sheet=LTW
$ws.Range("H27").Value = 975996.1
$ws.Range("I27").Value = 2533901.5
$ws.Range("K27").Value = 2533901.5
$ws.Range("M27").Value = -2533794.5

# hunk 31  sheet=LTW
$ws.Range("H136").Value = 20051400
$ws.Range("I136").Value = 3402737.2
$ws.Range("J136").Value = 66667660
$ws.Range("K136").Value = 10208211.6
$ws.Range("L136").Value = 200002980
$ws.Range("M136").Value = -10205661.6
$ws.Range("N136").Value = -200008080

# hunk 32  sheet=WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 26877.111
$ws.Range("J108").Value = 26877.111
$ws.Range("L108").Value = 26877.111
$ws.Range("N108").Value = -34557.111

# hunk 33  sheet=WVR
$ws.Range("H109").Value = 20340.334
$ws.Range("J109").Value = 20340.334
$ws.Range("L109").Value = 20340.334
$ws.Range("N109").Value = -23114.334

